$d = $word.ActiveDocument
for ($i = 0; $i -lt 6; $i++) {
  $cell = $d.Tables(1).Cell(2, 1)
  $full = $cell.Range
  $full.End = $full.End - 1
  $txt = $full.Text
  $full.Text = $txt
  Write-Host "iter $i text=[$txt]"
}
